# DataTable.xlsx edit
#
# Source-data correction: the "Aug 2025" release-date text in the data
# table (cell D13) is updated from a placeholder to the confirmed
# month, matching the commit message's "Correct data source text".
#
# The sheet's saved selection is also moved onto D14 (the last data
# row), mirroring where the author left the cursor when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Aug 2025" data-source/period text.
$ws.Range("D13").Value = "Aug 2025 (Sep 25)"

# Leave the selection on D14, as in the saved workbook.
$ws.Range("D14").Select() | Out-Null
